# Daily attendance processing - 2025-11-20 05:22:54
#
# The "Recorded By" column (G) contains comma-separated lists of users who
# recorded/touched each attendance entry (e.g. "someone@example.com, System").
# This pass re-orders each list so that the literal value "System" (exact
# case) is swapped from wherever it sits to the front of the list, swapping
# places with whichever entry is currently first. Lists that do not contain
# the literal "System" value are left untouched.
#
# Note: string comparisons in this runtime are case-INsensitive by default
# (-eq/-ne/-contains), so exact-case matching/inequality checks below use the
# .NET [string]::Equals method, which is case-sensitive.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1

# Header is row 1 ("Recorded By" lives in column G / index 7); data starts row 2.
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $current = $cell.Value2

    if ($null -eq $current) { continue }
    if ($current -eq "") { continue }

    $parts = @($current -split ", ")

    if ($parts.Length -lt 2) { continue }

    $hasExactSystem = $false
    foreach ($p in $parts) {
        if ($p.Equals("System")) {
            $hasExactSystem = $true
        }
    }

    if ($hasExactSystem) {
        $firstPart = $parts[0]
        $lastIndex = $parts.Length - 1
        $lastPart = $parts[$lastIndex]

        $parts[0] = $lastPart
        $parts[$lastIndex] = $firstPart

        $updated = ($parts -join ", ")

        if (-not $updated.Equals($current)) {
            $cell.Value = $updated
        }
    }
}
